$wb = $excel.ActiveWorkbook

# --- Sheet1: shrink autofilter range from A1:XFD to A1 ---
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.AutoFilter.Range.Value | Out-Null
$ws1.Range("A1").AutoFilter()

# --- Sheet2: timeline update for "LID Board" test fixture dev ---
$ws2 = $wb.Worksheets.Item("Sheet2")

# Change the section title to a new string "TF Dev LID Board"
$ws2.Range("D18").Value = "TF Dev LID Board"

# Update the active selection to J21
$ws2.Range("J21").Select() | Out-Null

# Update progress % and completion date values in the new timeline block
$ws2.Range("F20").Value = 0
$ws2.Range("H20").Value = 43630

$ws2.Range("F21").Value = 0
$ws2.Range("H21").Value = 43633

$ws2.Range("F22").Value = 0
$ws2.Range("H22").Value = 43634

$ws2.Range("H23").Value = 43636

$ws2.Range("H24").Value = 43637

$ws2.Range("F25").Value = 0
$ws2.Range("H25").Value = 43641

$ws2.Range("H26").Value = 43642

$ws2.Range("H27").Value = 43642

$ws2.Range("H28").Value = 43644

$wb.Save()
